# side_by_side_tests.xlsx — add two more "Raw/Provided Scientific Name" search
# tests, inserted above the existing mapping-legend/static-map rows, and leave
# the selection where the author's cursor ended up after typing them in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new test rows just above the old row 4 (old rows 4-5
# shift down to 6-7). Excel's Insert() carries the formatting of the row
# above down into the new rows, so B4/B5 already pick up the commented-value
# style used by B1:B3.
$ws.Rows("4:5").Insert()

# Row 5 first: "Acacia dealbata" test (typed in before the Osphranter one,
# so its strings land earlier in the shared-string table).
$ws.Range("A5").Value = "Search for Raw/Provided Scientific Name 'Acacia dealbata' should turn up assorted silver wattles"
$ws.Range("C5").Value = "occurrences/search"
$ws.Range("D5").Value = "q=raw_name%3A%22Acacia%20dealbata%22&start=0&pageSize=20&sort=first_loaded_date&dir=desc&qc=&facets=taxon_name"

# Row 4: "Osphranter rufus" test.
$ws.Range("A4").Value = "Search for Raw/Provided Scientific Name 'Osphranter rufus' should turn up assorted red kangaroos"
$ws.Range("C4").Value = "occurrences/search"
$ws.Range("D4").Value = "q=raw_name%3A%22Osphranter%20rufus%22&start=0&pageSize=20&sort=first_loaded_date&dir=desc&qc=&facets=taxon_name"

# Give the A/D cells of both new rows the same "commented (hardcoded) value"
# look already used for the empty B1:B3 cells (Monaco, 9pt, greenish text) —
# copy its format instead of poking Font.* piecemeal so we reuse the
# existing style instead of minting new font/xf entries.
$ws.Range("B4").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)

# Leave the selection where the author's cursor ended up next (just past the
# bottom-right of the new data).
$ws.Range("D8").Select()
